# Wind Constant in test script lowered. Checked each test for stability.
#
# This edit:
#   1. Adds a new "Stable" flag column (P) to the flag_test_data sheet,
#      classifying each of the 121 test rows as Stable / unstable / (s)
#      based on the (lowered) wind-constant stability re-check.
#   2. Re-selects the flag_test_data sheet (was "Explicit Euler") and
#      updates its scroll/selection state to where the review left off.

$wb = $excel.ActiveWorkbook

# --- flag_test_data sheet: populate the new "Stable" column (P1:P121) ---
$ws1 = $wb.Worksheets.Item("flag_test_data")

$arrP = New-Object 'object[,]' 121,1
$arrP[0,0] = "Stable"
$arrP[1,0] = "Stable"
$arrP[2,0] = "Stable"
$arrP[3,0] = "u"
$arrP[4,0] = "u"
$arrP[5,0] = "u"
$arrP[6,0] = "Stable"
$arrP[7,0] = "Stable"
$arrP[8,0] = "u"
$arrP[9,0] = "u"
$arrP[10,0] = "u"
$arrP[11,0] = "Stable"
$arrP[12,0] = "u"
$arrP[13,0] = "u"
$arrP[14,0] = "u"
$arrP[15,0] = "u"
$arrP[16,0] = "Stable"
$arrP[17,0] = "u"
$arrP[18,0] = "u"
$arrP[19,0] = "u"
$arrP[20,0] = "u"
$arrP[21,0] = "Stable"
$arrP[22,0] = "u"
$arrP[23,0] = "u"
$arrP[24,0] = "u"
$arrP[25,0] = "u"
$arrP[26,0] = "Stable"
$arrP[27,0] = "u"
$arrP[28,0] = "u"
$arrP[29,0] = "u"
$arrP[30,0] = "u"
$arrP[31,0] = "Stable"
$arrP[32,0] = "Stable"
$arrP[33,0] = "Stable"
$arrP[34,0] = "Stable"
$arrP[35,0] = "Stable"
$arrP[36,0] = "Stable"
$arrP[37,0] = "Stable"
$arrP[38,0] = "Stable"
$arrP[39,0] = "u"
$arrP[40,0] = "u"
$arrP[41,0] = "Stable"
$arrP[42,0] = "Stable"
$arrP[43,0] = "u"
$arrP[44,0] = "u"
$arrP[45,0] = "u"
$arrP[46,0] = "Stable"
$arrP[47,0] = "Stable"
$arrP[48,0] = "u"
$arrP[49,0] = "u"
$arrP[50,0] = "u"
$arrP[51,0] = "Stable"
$arrP[52,0] = "u"
$arrP[53,0] = "u"
$arrP[54,0] = "u"
$arrP[55,0] = "u"
$arrP[56,0] = "Stable"
$arrP[57,0] = "u"
$arrP[58,0] = "u"
$arrP[59,0] = "u"
$arrP[60,0] = "u"
$arrP[61,0] = "Stable"
$arrP[62,0] = "Stable"
$arrP[63,0] = "Stable"
$arrP[64,0] = "Stable"
$arrP[65,0] = "u"
$arrP[66,0] = "Stable"
$arrP[67,0] = "Stable"
$arrP[68,0] = "u"
$arrP[69,0] = "u"
$arrP[70,0] = "u"
$arrP[71,0] = "Stable"
$arrP[72,0] = "Stable"
$arrP[73,0] = "u"
$arrP[74,0] = "u"
$arrP[75,0] = "u"
$arrP[76,0] = "Stable"
$arrP[77,0] = "u"
$arrP[78,0] = "u"
$arrP[79,0] = "u"
$arrP[80,0] = "u"
$arrP[81,0] = "Stable"
$arrP[82,0] = "u"
$arrP[83,0] = "u"
$arrP[84,0] = "u"
$arrP[85,0] = "u"
$arrP[86,0] = "Stable"
$arrP[87,0] = "u"
$arrP[88,0] = "u"
$arrP[89,0] = "u"
$arrP[90,0] = "u"
$arrP[91,0] = "s"
$arrP[92,0] = "s"
$arrP[93,0] = "s"
$arrP[94,0] = "u"
$arrP[95,0] = "u"
$arrP[96,0] = "s"
$arrP[97,0] = "u"
$arrP[98,0] = "u"
$arrP[99,0] = "u"
$arrP[100,0] = "u"
$arrP[101,0] = "s"
$arrP[102,0] = "u"
$arrP[103,0] = "u"
$arrP[104,0] = "u"
$arrP[105,0] = "u"
$arrP[106,0] = "s"
$arrP[107,0] = "u"
$arrP[108,0] = "u"
$arrP[109,0] = "u"
$arrP[110,0] = "u"
$arrP[111,0] = "Stable"
$arrP[112,0] = "u"
$arrP[113,0] = "u"
$arrP[114,0] = "u"
$arrP[115,0] = "u"
$arrP[116,0] = "s"
$arrP[117,0] = "u"
$arrP[118,0] = "u"
$arrP[119,0] = "u"
$arrP[120,0] = "u"
$ws1.Range("P1:P121").Value = $arrP

# --- flag_test_data becomes the active / selected sheet again ---
$ws1.Activate() | Out-Null
$ws1.Range("T108").Select() | Out-Null

# --- Explicit Euler (previously active) is no longer the selected tab ---
$ws5 = $wb.Worksheets.Item("Explicit Euler")
